# Adds "Stored in refrigerator" as a new option for the 'preparation_condition'
# lookup list, wires it into the Sample Section sheet's data validation range,
# and bumps the .metadata sheet's pav:createdOn timestamp.

$wb = $excel.ActiveWorkbook

# 1. Append the new lookup row to the 'preparation_condition' sheet
#    (was A1:B7, becomes A1:B8).
$prepCond = $wb.Worksheets.Item("preparation_condition")
$prepCond.Range("A8").Value = "Stored in refrigerator"
$prepCond.Range("B8").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000104"

# 2. Extend the Sample Section sheet's column H (preparation_condition)
#    validation list so it covers the newly added row.
$sampleSection = $wb.Worksheets.Item("Sample Section")
$sampleSection.Range("H2:H1001").Validation.Formula1 = '=''preparation_condition''!$A$1:$A$8'

# 3. Record the new pav:createdOn timestamp on the .metadata sheet.
$metadata = $wb.Worksheets.Item(".metadata")
$metadata.Range("C2").Value = "2024-02-21T09:25:13-08:00"
